$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.619.06'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '3.180.22'
$ws.Range("E3").Value = '  -3.22%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.19'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '614.86'
$ws.Range("E6").Value = '  -3.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.388'
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.688'
$ws.Range("E8").Value = '  -6.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D10").Value = '3.175.78'
$ws.Range("E10").Value = '  -3.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.574'
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.176'
$ws.Range("E12").Value = '  -5.35%  '
$ws.Range("E13").Value = '  -6.13%  '
$ws.Range("D14").Value = '3.774.41'
$ws.Range("E14").Value = '  -2.79%  '
$ws.Range("D15").Value = '89.584.08'
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.76'
$ws.Range("E16").Value = '  -6.15%  '
$ws.Range("E17").Value = '  -5.59%  '
$ws.Range("D18").Value = '3.179.27'
$ws.Range("E18").Value = '  -3.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.26'
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.36'
$ws.Range("E20").Value = '  -6.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '434.68'
$ws.Range("E21").Value = '  -1.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000197'
$ws.Range("E22").Value = '  +35.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.56'
$ws.Range("E23").Value = '  -5.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.05'
$ws.Range("E24").Value = '  -6.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.11'
$ws.Range("E25").Value = '  -4.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.66'
$ws.Range("E26").Value = '  -6.27%  '
$ws.Range("D27").Value = '3.349.71'
$ws.Range("E27").Value = '  -2.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '75.17'
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  -7.21%  '
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.04'
$ws.Range("E32").Value = '  +26.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.39'
$ws.Range("E33").Value = '  -6.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '532.86'
$ws.Range("E34").Value = '  -7.31%  '
$ws.Range("E35").Value = '  -5.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.85'
$ws.Range("E36").Value = '  -6.75%  '
$ws.Range("E37").Value = '  -9.44%  '
$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.32'
$ws.Range("E38").Value = '  +2.16%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.89'
$ws.Range("E39").Value = '  -5.71%  '
$ws.Range("E40").Value = '  -9.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.91'
$ws.Range("E43").Value = '  -7.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.370'
$ws.Range("E44").Value = '  -9.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '149.09'
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.59'
$ws.Range("E46").Value = '  -2.86%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '172.12'
$ws.Range("E47").Value = '  -5.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.123'
$ws.Range("E48").Value = '  -9.81%  '
$ws.Range("E49").Value = '  -9.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.05'
$ws.Range("E50").Value = '  -5.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.607'
$ws.Range("E51").Value = '  -4.83%  '
